# Scheduled market-data refresh: update cached price/profit figures
# across the Leve profit tables (Table_<Job>) on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39947.5
$ws.Range("J3").Value = 39947.5
$ws.Range("L3").Value = 39947.5
$ws.Range("N3").Value = -40175.5
$ws.Range("H20").Value = 26000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 26000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 26000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -26460
$ws.Range("H32").Value = 27128816
$ws.Range("I32").Value = 83333600
$ws.Range("J32").Value = 4646902.5
$ws.Range("K32").Value = 83333600
$ws.Range("L32").Value = 4646902.5
$ws.Range("M32").Value = -83333274
$ws.Range("N32").Value = -4647554.5
$ws.Range("H34").Value = 10033
$ws.Range("I34").Value = 1480.5714
$ws.Range("J34").Value = 24999.75
$ws.Range("K34").Value = 1480.5714
$ws.Range("L34").Value = 24999.75
$ws.Range("M34").Value = -1277.5714
$ws.Range("N34").Value = -25405.75
$ws.Range("H35").Value = 26000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 26000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 26000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -26758
$ws.Range("H36").Value = 10033
$ws.Range("I36").Value = 1480.5714
$ws.Range("J36").Value = 24999.75
$ws.Range("K36").Value = 1480.5714
$ws.Range("L36").Value = 24999.75
$ws.Range("M36").Value = -765.5714
$ws.Range("N36").Value = -26429.75
$ws.Range("H47").Value = 21000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 21000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 21000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -22944
$ws.Range("H51").Value = 8538.538
$ws.Range("I51").Value = 8300.333000000001
$ws.Range("K51").Value = 8300.333000000001
$ws.Range("M51").Value = -7816.333000000001
$ws.Range("H93").Value = 85233.5
$ws.Range("J93").Value = 85233.5
$ws.Range("L93").Value = 85233.5
$ws.Range("N93").Value = -90225.5
$ws.Range("H95").Value = 20624
$ws.Range("J95").Value = 20624
$ws.Range("L95").Value = 20624
$ws.Range("N95").Value = -26116
$ws.Range("H102").Value = 39947.5
$ws.Range("J102").Value = 39947.5
$ws.Range("L102").Value = 39947.5
$ws.Range("N102").Value = -46437.5
$ws.Range("H105").Value = 39500
$ws.Range("J105").Value = 39500
$ws.Range("L105").Value = 39500
$ws.Range("N105").Value = -46488
$ws.Range("H116").Value = 39598.31
$ws.Range("I116").Value = 52998.668
$ws.Range("J116").Value = 4422.375
$ws.Range("K116").Value = 52998.668
$ws.Range("L116").Value = 4422.375
$ws.Range("M116").Value = -49556.668
$ws.Range("N116").Value = -11306.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14177.071
$ws.Range("I28").Value = 8196.200000000001
$ws.Range("J28").Value = 29129.25
$ws.Range("K28").Value = 8196.200000000001
$ws.Range("L28").Value = 29129.25
$ws.Range("M28").Value = -8004.200000000001
$ws.Range("N28").Value = -29513.25
$ws.Range("H31").Value = 6800
$ws.Range("I31").Value = 6800
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6800
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6506
$ws.Range("N31").ClearContents()
$ws.Range("H93").Value = 27000
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H99").Value = 14177.071
$ws.Range("I99").Value = 8196.200000000001
$ws.Range("J99").Value = 29129.25
$ws.Range("K99").Value = 8196.200000000001
$ws.Range("L99").Value = 29129.25
$ws.Range("M99").Value = -5201.200000000001
$ws.Range("N99").Value = -35119.25
$ws.Range("H132").Value = 33541.344
$ws.Range("I132").Value = 51532.4
$ws.Range("J132").Value = 3556.25
$ws.Range("K132").Value = 154597.2
$ws.Range("L132").Value = 10668.75
$ws.Range("M132").Value = -152067.2
$ws.Range("N132").Value = -15728.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 28189
$ws.Range("J21").Value = 28189
$ws.Range("L21").Value = 28189
$ws.Range("N21").Value = -28661
$ws.Range("H28").Value = 29709
$ws.Range("J28").Value = 29709
$ws.Range("L28").Value = 29709
$ws.Range("N28").Value = -30297
$ws.Range("H38").Value = 11400
$ws.Range("J38").Value = 11400
$ws.Range("L38").Value = 11400
$ws.Range("N38").Value = -12232
$ws.Range("H44").Value = 21050
$ws.Range("J44").Value = 21050
$ws.Range("L44").Value = 21050
$ws.Range("N44").Value = -22044
$ws.Range("H101").Value = 22000
$ws.Range("J101").Value = 22000
$ws.Range("L101").Value = 22000
$ws.Range("N101").Value = -28490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 83335400
$ws.Range("I35").Value = 125000610
$ws.Range("K35").Value = 125000610
$ws.Range("M35").Value = -125000316
$ws.Range("H38").Value = 5324.95
$ws.Range("I38").Value = 499.66666
$ws.Range("J38").Value = 6176.4707
$ws.Range("K38").Value = 499.66666
$ws.Range("L38").Value = 6176.4707
$ws.Range("M38").Value = -122.66666
$ws.Range("N38").Value = -6930.4707
$ws.Range("H46").Value = 5324.95
$ws.Range("I46").Value = 499.66666
$ws.Range("J46").Value = 6176.4707
$ws.Range("K46").Value = 499.66666
$ws.Range("L46").Value = 6176.4707
$ws.Range("M46").Value = -288.66666
$ws.Range("N46").Value = -6598.4707
$ws.Range("H96").Value = 10082.385
$ws.Range("J96").Value = 10082.385
$ws.Range("L96").Value = 10082.385
$ws.Range("N96").Value = -15574.385
$ws.Range("H107").Value = 500.5
$ws.Range("I107").Value = 353.2
$ws.Range("J107").Value = 701.36365
$ws.Range("K107").Value = 353.2
$ws.Range("L107").Value = 701.36365
$ws.Range("M107").Value = 1566.8
$ws.Range("N107").Value = -4541.36365
$ws.Range("H112").Value = 35000
$ws.Range("J112").Value = 35000
$ws.Range("L112").Value = 35000
$ws.Range("N112").Value = -37954

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2010.3334
$ws.Range("I31").Value = 731
$ws.Range("J31").Value = 2650
$ws.Range("K31").Value = 731
$ws.Range("L31").Value = 2650
$ws.Range("M31").Value = -439
$ws.Range("N31").Value = -3234
$ws.Range("H37").Value = 2010.3334
$ws.Range("I37").Value = 731
$ws.Range("J37").Value = 2650
$ws.Range("K37").Value = 731
$ws.Range("L37").Value = 2650
$ws.Range("M37").Value = -454
$ws.Range("N37").Value = -3204
$ws.Range("H94").Value = 27000
$ws.Range("J94").Value = 27000
$ws.Range("L94").Value = 27000
$ws.Range("N94").Value = -28352
$ws.Range("H98").Value = 20544.2
$ws.Range("J98").Value = 20544.2
$ws.Range("L98").Value = 20544.2
$ws.Range("N98").Value = -26534.2
$ws.Range("H99").Value = 6081.7144
$ws.Range("I99").Value = 3762
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 3762
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -1516
$ws.Range("N99").Value = -24492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 434.66666
$ws.Range("I9").Value = 321.6
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 321.6
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -97.60000000000002
$ws.Range("N9").Value = -1448
$ws.Range("H29").Value = 5933.3335
$ws.Range("I29").Value = 3900
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 3900
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -3605
$ws.Range("N29").Value = -10590
$ws.Range("H57").Value = 9420.143
$ws.Range("J57").Value = 15225
$ws.Range("L57").Value = 15225
$ws.Range("N57").Value = -16357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 5886.2
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 5590.4
$ws.Range("I34").Value = 5476
$ws.Range("J34").Value = 5666.6665
$ws.Range("K34").Value = 5476
$ws.Range("L34").Value = 5666.6665
$ws.Range("M34").Value = -5273
$ws.Range("N34").Value = -6072.6665
$ws.Range("H61").Value = 8238.5
$ws.Range("I61").Value = 3850.3333
$ws.Range("J61").Value = 12626.667
$ws.Range("K61").Value = 3850.3333
$ws.Range("L61").Value = 12626.667
$ws.Range("M61").Value = -3558.3333
$ws.Range("N61").Value = -13210.667
$ws.Range("H132").Value = 1767.2075
$ws.Range("I132").Value = 1212.6052
$ws.Range("J132").Value = 3172.2
$ws.Range("K132").Value = 3637.8156
$ws.Range("L132").Value = 9516.599999999999
$ws.Range("M132").Value = -1107.8156
$ws.Range("N132").Value = -14576.6
$ws.Range("H136").Value = 17811272
$ws.Range("I136").Value = 21064000
$ws.Range("J136").Value = 7402537
$ws.Range("K136").Value = 63192000
$ws.Range("L136").Value = 22207611
$ws.Range("M136").Value = -63189450
$ws.Range("N136").Value = -22212711

